# Generate Report for Handback
# Applies the localization-status.xlsx handback update:
#  - Overview sheet: rows for "file_partial_deleted_multi_path_1.md" and
#    "file_no_deleted_single_path.md" swap places (row 4 / row 5), and the
#    three "no_deleted_multi_path_*" / "partial_deleted" rows move from
#    "Ready for handoff" to "Handed back: in sync with en-US".
#  - zh-cn / de-de sheets: same row swap, status updates, and the
#    "Latest Target File" / "Latest Handback File" / "Latest Handback
#    DateTime" columns get populated for the handed-back rows, with an
#    out-of-date warning recorded for the single-path file.

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet: Overview
# -------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 2 - file_no_deleted_multi_path_1.md : status moves to Handed back
$ws.Range("E2").Value = "Handed back: in sync with en-US"
$ws.Range("F2").Value = "Handed back: in sync with en-US"

# Row 3 - file_no_deleted_multi_path_2.md : status moves to Handed back
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"

# Row 4 / Row 5 swap identity: partial_deleted moves to row 4,
# no_deleted_single_path moves to row 5.
$ws.Range("A4").Value = "file_partial_deleted_multi_path_1.md"
$ws.Range("B4").Value = "e2e\file_partial_deleted_multi_path_1.md"
$ws.Range("E4").Value = "Handed back: in sync with en-US"
$ws.Range("F4").Value = "Handed back: in sync with en-US"
$ws.Range("G4").Value = "2016-08-30 14:45:40"

$ws.Range("A5").Value = "file_no_deleted_single_path.md"
$ws.Range("B5").Value = "e2e\file_no_deleted_single_path.md"
$ws.Range("E5").Value = "Ready for handoff"
$ws.Range("F5").Value = "Ready for handoff"
$ws.Range("G5").Value = "2016-08-30 14:46:15"

# Rebuild the B2:B6 hyperlinks, preserving each relationship's original
# target address but refreshing the display text to match the new row
# contents (row 4 / row 5 display text is swapped, matching the other
# rows keeping their original target).
$links = @()
foreach ($h in $ws.Hyperlinks) { $links += $h }
$links[0].TextToDisplay = "e2e\file_no_deleted_multi_path_1.md"
$links[1].TextToDisplay = "e2e\file_no_deleted_multi_path_2.md"
$links[2].TextToDisplay = "e2e\file_partial_deleted_multi_path_1.md"
$links[3].TextToDisplay = "e2e\file_no_deleted_single_path.md"
$links[4].TextToDisplay = "e2e\newfile.8fcde224-5f0e-49fe-bd42-5ad52d8ea82a.md"

# -------------------------------------------------------------------------
# Helper data shared by the zh-cn / de-de sheets
# -------------------------------------------------------------------------
$urlMultiPath1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee94f4816613fe36c27c95a2d490c9f47944ddf3/e2e/file_no_deleted_multi_path_1.md"
$urlMultiPath2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a69ed39e2504b0f277c93f6ab6084297bc0ef6c/e2e/file_no_deleted_multi_path_2.md"
$urlSinglePath = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a69ed39e2504b0f277c93f6ab6084297bc0ef6c/e2e/file_no_deleted_single_path.md"
$urlPartialDeleted = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee94f4816613fe36c27c95a2d490c9f47944ddf3/e2e/file_partial_deleted_multi_path_1.md"

function Update-LangSheet([string]$sheetName, [string]$lang, [string]$handbackTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2 - file_no_deleted_multi_path_1.md
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("I2").Value = "file_no_deleted_multi_path_1.md"
    $ws.Range("J2").Value = ("file_no_deleted_multi_path_1.2336e4b28de82563dc9ea89a1eb254ab730456de." + $lang + ".xlf")
    $ws.Range("K2").Value = $handbackTime

    # Row 3 - file_no_deleted_multi_path_2.md (content duplicate of row 2)
    $ws.Range("C3").Value = "Handed back: in sync with en-US"
    $ws.Range("I3").Value = "file_no_deleted_multi_path_1.md"
    $ws.Range("J3").Value = ("file_no_deleted_multi_path_1.2336e4b28de82563dc9ea89a1eb254ab730456de." + $lang + ".xlf")
    $ws.Range("K3").Value = $handbackTime

    # Row 4 / Row 5 swap identity: partial_deleted moves to row 4,
    # no_deleted_single_path moves to row 5.
    $ws.Range("A4").Value = "file_partial_deleted_multi_path_1.md"
    $ws.Range("C4").Value = "Handed back: in sync with en-US"
    $ws.Range("G4").Value = ("file_partial_deleted_multi_path_1.950eb2cafff7c1eccbf1f4ff641b1cbe66aecd0c." + $lang + ".xlf")
    $ws.Range("H4").Value = "2016-08-30 14:45:30"
    $ws.Range("I4").Value = "file_partial_deleted_multi_path_1.md"
    $ws.Range("J4").Value = ("file_partial_deleted_multi_path_1.950eb2cafff7c1eccbf1f4ff641b1cbe66aecd0c." + $lang + ".xlf")
    $ws.Range("K4").Value = $handbackTime

    $ws.Range("A5").Value = "file_no_deleted_single_path.md"
    $ws.Range("G5").Value = ("file_no_deleted_single_path.84480e30a75eaa7877ca77d49d2059db6121edd7." + $lang + ".xlf")
    $ws.Range("H5").Value = "2016-08-30 14:45:57"
    $ws.Range("I5").Value = "file_no_deleted_single_path.md"
    $ws.Range("J5").Value = ("file_no_deleted_single_path.84480e30a75eaa7877ca77d49d2059db6121edd7." + $lang + ".xlf")
    $ws.Range("K5").Value = $handbackTime
    $ws.Range("P5").Value = ("The version of handback file is not the latest, current: " + $urlPartialDeleted + ", latest: " + $urlSinglePath + ".")

    # Rebuild hyperlinks: column A keeps one hyperlink per row (2-6),
    # column I gets a new hyperlink for rows 2-5 pointing at the same
    # source file as that row's "Latest Target File" value.
    $ws.Range("A1").Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $urlMultiPath1, "", "", "file_no_deleted_multi_path_1.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlMultiPath1, "", "", "file_no_deleted_multi_path_1.md") | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), $urlMultiPath2, "", "", "file_no_deleted_multi_path_2.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlMultiPath1, "", "", "file_no_deleted_multi_path_1.md") | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A4"), $urlPartialDeleted, "", "", "file_partial_deleted_multi_path_1.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I4"), $urlPartialDeleted, "", "", "file_partial_deleted_multi_path_1.md") | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A5"), $urlSinglePath, "", "", "file_no_deleted_single_path.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I5"), $urlSinglePath, "", "", "file_no_deleted_single_path.md") | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a69ed39e2504b0f277c93f6ab6084297bc0ef6c/e2e/newfile.8fcde224-5f0e-49fe-bd42-5ad52d8ea82a.md", "", "", "newfile.8fcde224-5f0e-49fe-bd42-5ad52d8ea82a.md") | Out-Null

    # Column widths widened for the new, longer content in columns C/I/J/P.
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 34.9774257114955
    $ws.Columns.Item(10).ColumnWidth = 40
    $ws.Columns.Item(16).ColumnWidth = 40
}

Update-LangSheet "zh-cn" "zh-cn" "2016-08-30 14:46:36"
Update-LangSheet "de-de" "de-de" "2016-08-30 14:46:43"

# Overview sheet column widths for the zh-cn/de-de status columns (now hold
# the longer "Handed back: in sync with en-US" text).
$ws = $wb.Worksheets.Item("Overview")
$ws.Columns.Item(5).ColumnWidth = 29.9777047293527
$ws.Columns.Item(6).ColumnWidth = 29.9777047293527
